# "writeStatement" (row 50) completion went from 0% to 20%.
# C50 (=IF(B50>0,1,0)) and the header rollups B1/C1 (=SUM(..)/49) are
# formulas, so they recompute automatically once B50 changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B50").Value = 0.2

# Make sure every dependent formula (C50, B1, C1) is recalculated.
$excel.Calculate()

# The view scrolled down (row 21 became the top visible row) and the
# selection moved from E36 to C50.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C50").Select()
